$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 26600
$ws.Range("J3").Value = 26600
$ws.Range("L3").Value = 26600
$ws.Range("N3").Value = -26828
$ws.Range("H20").Value = 1299.6666
$ws.Range("I20").Value = 1299.6666
$ws.Range("K20").Value = 1299.6666
$ws.Range("M20").Value = -1069.6666
$ws.Range("H33").Value = 475.66666
$ws.Range("I33").Value = 550
$ws.Range("J33").Value = 438.5
$ws.Range("K33").Value = 550
$ws.Range("L33").Value = 438.5
$ws.Range("M33").Value = -321
$ws.Range("N33").Value = -896.5
$ws.Range("H35").Value = 1299.6666
$ws.Range("I35").Value = 1299.6666
$ws.Range("K35").Value = 1299.6666
$ws.Range("M35").Value = -920.6666
$ws.Range("H88").Value = 1326.5714
$ws.Range("J88").Value = 1557.2
$ws.Range("L88").Value = 1557.2
$ws.Range("N88").Value = -2369.2
$ws.Range("H91").Value = 1326.5714
$ws.Range("J91").Value = 1557.2
$ws.Range("L91").Value = 1557.2
$ws.Range("N91").Value = -4365.2
$ws.Range("H102").Value = 26600
$ws.Range("J102").Value = 26600
$ws.Range("L102").Value = 26600
$ws.Range("N102").Value = -33090
$ws.Range("H138").Value = 3097.8118
$ws.Range("J138").Value = 3023.1904
$ws.Range("L138").Value = 9069.5712
$ws.Range("N138").Value = -19349.5712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2995.2354
$ws.Range("I32").Value = 3637.25
$ws.Range("K32").Value = 3637.25
$ws.Range("M32").Value = -3350.25
$ws.Range("H44").Value = 69998
$ws.Range("J44").Value = 69998
$ws.Range("L44").Value = 69998
$ws.Range("N44").Value = -70974
$ws.Range("H55").Value = 48723.5
$ws.Range("J55").Value = 69998
$ws.Range("L55").Value = 69998
$ws.Range("N55").Value = -70628
$ws.Range("H110").Value = 1090.6
$ws.Range("I110").Value = 1119.5
$ws.Range("J110").Value = 975
$ws.Range("K110").Value = 1119.5
$ws.Range("L110").Value = 975
$ws.Range("M110").Value = 925.5
$ws.Range("N110").Value = -5065

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 46656.168
$ws.Range("I82").Value = 6665
$ws.Range("K82").Value = 6665
$ws.Range("M82").Value = -6282
$ws.Range("H85").Value = 46656.168
$ws.Range("I85").Value = 6665
$ws.Range("K85").Value = 6665
$ws.Range("M85").Value = -5339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3654.6667
$ws.Range("I16").Value = 1400
$ws.Range("J16").Value = 4298.857
$ws.Range("K16").Value = 1400
$ws.Range("L16").Value = 4298.857
$ws.Range("M16").Value = -1113
$ws.Range("N16").Value = -4872.857
$ws.Range("H25").Value = 2000
$ws.Range("J25").Value = 2000
$ws.Range("L25").Value = 2000
$ws.Range("N25").Value = -2348
$ws.Range("H50").Value = 59998.5
$ws.Range("J50").Value = 59998.5
$ws.Range("L50").Value = 59998.5
$ws.Range("N50").Value = -61248.5
$ws.Range("H58").Value = 3950
$ws.Range("I58").Value = 3950
$ws.Range("K58").Value = 3950
$ws.Range("M58").Value = -3747
$ws.Range("H62").Value = 6998.5
$ws.Range("J62").Value = 8333
$ws.Range("L62").Value = 8333
$ws.Range("N62").Value = -9581
$ws.Range("H65").Value = 6998.5
$ws.Range("J65").Value = 8333
$ws.Range("L65").Value = 41665
$ws.Range("N65").Value = -47905
$ws.Range("H113").Value = 3654.6667
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 4298.857
$ws.Range("K113").Value = 1400
$ws.Range("L113").Value = 4298.857
$ws.Range("M113").Value = 770
$ws.Range("N113").Value = -8638.857
$ws.Range("H132").Value = 1371.2
$ws.Range("I132").Value = 1403.2
$ws.Range("K132").Value = 4209.6
$ws.Range("M132").Value = -1679.6
$ws.Range("H134").Value = 2484.2
$ws.Range("I134").Value = 2484.2
$ws.Range("K134").Value = 7452.599999999999
$ws.Range("M134").Value = -4917.599999999999
$ws.Range("H136").Value = 3950
$ws.Range("I136").Value = 3950
$ws.Range("K136").Value = 11850
$ws.Range("M136").Value = -9300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3792.8333
$ws.Range("I3").Value = 3792.8333
$ws.Range("K3").Value = 11378.4999
$ws.Range("M3").Value = -11266.4999
$ws.Range("H48").Value = 19666.666
$ws.Range("J48").Value = 19666.666
$ws.Range("L48").Value = 58999.99800000001
$ws.Range("N48").Value = -59499.99800000001
$ws.Range("H56").Value = 4788.981
$ws.Range("I56").Value = 4788.981
$ws.Range("K56").Value = 4788.981
$ws.Range("M56").Value = -4258.981
$ws.Range("H108").Value = 3181.923
$ws.Range("I108").Value = 1033.2727
$ws.Range("K108").Value = 3099.8181
$ws.Range("M108").Value = -219.8181
$ws.Range("H121").Value = 68975.5
$ws.Range("I121").Value = 25519.75
$ws.Range("J121").Value = 79839.44
$ws.Range("K121").Value = 76559.25
$ws.Range("L121").Value = 239518.32
$ws.Range("M121").Value = -75249.25
$ws.Range("N121").Value = -242138.32
$ws.Range("H124").Value = 19111
$ws.Range("I124").Value = 4000
$ws.Range("J124").Value = 22888.75
$ws.Range("K124").Value = 12000
$ws.Range("L124").Value = 68666.25
$ws.Range("M124").Value = -7090
$ws.Range("N124").Value = -78486.25
$ws.Range("H128").Value = 89666.664
$ws.Range("I128").Value = 89666.664
$ws.Range("K128").Value = 268999.992
$ws.Range("M128").Value = -264019.992
$ws.Range("H129").Value = 3719.4666
$ws.Range("I129").Value = 2099
$ws.Range("K129").Value = 6297
$ws.Range("M129").Value = -1297

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 60000
$ws.Range("J103").Value = 60000
$ws.Range("L103").Value = 60000
$ws.Range("N103").Value = -62344
$ws.Range("H113").Value = 1881.6428
$ws.Range("I113").Value = 1436.6
$ws.Range("J113").Value = 2994.25
$ws.Range("K113").Value = 1436.6
$ws.Range("L113").Value = 2994.25
$ws.Range("M113").Value = 733.4000000000001
$ws.Range("N113").Value = -7334.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3275.8
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H68").Value = 2062.2222
$ws.Range("I68").Value = 2062.2222
$ws.Range("K68").Value = 2062.2222
$ws.Range("M68").Value = -1313.2222
$ws.Range("H71").Value = 2062.2222
$ws.Range("I71").Value = 2062.2222
$ws.Range("K71").Value = 10311.111
$ws.Range("M71").Value = -6567.111000000001
$ws.Range("H93").Value = 2499.5
$ws.Range("I93").Value = 2499.5
$ws.Range("K93").Value = 2499.5
$ws.Range("M93").Value = -1251.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 50358.668
$ws.Range("I53").Value = 55038
$ws.Range("K53").Value = 55038
$ws.Range("M53").Value = -54431
